$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("ID", "Фамилия", "Имя", "Пол", "Рейтинг КОФНТ", "Рейтинг ФНТР", "Рейтинг РТТФ", "Дата рождения", "Разряд", "Город")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data row
$ws.Cells.Item(2, 1).Value = "azaza"
$ws.Cells.Item(2, 2).Value = "Богданов"
$ws.Cells.Item(2, 3).Value = "Роман"
$ws.Cells.Item(2, 4).Value = "М"
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 6
$ws.Cells.Item(2, 7).Value = 7
$ws.Cells.Item(2, 8).Value = 36526
$ws.Cells.Item(2, 8).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2, 9).Value = "МС"
$ws.Cells.Item(2, 10).Value = "Калининград"

# Column widths (autofit to match bestFit columns in diff)
$ws.Range("A1:J2").EntireColumn.AutoFit()

$ws.Range("L5").Select()
